$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G header text: "Value list" -> "Value"
$ws.Range("G3").Value = "Value"

# Expand the "Entity Template" header cell from A1:B1 down to A1:B2
# and give the merged block top-aligned wrapped text.
$ws.Range("A1").VerticalAlignment = -4160
$ws.Range("A1").WrapText = $true
$ws.Range("A1:B1").MergeCells = $false
$ws.Range("A1:B2").MergeCells = $true

# Move the saved selection to E4
$ws.Range("E4").Select()

Write-Output "edit applied"
